# Automatic map refresh: the data rows (A2:P51) get re-ordered because the
# upstream source re-emits the 5 most-recent records interleaved into their
# chronological slot among the already-known records. Every column travels
# together with its row, so this is a pure row permutation - no cell text
# needs to be retyped (and thus no re-encoding risk for the Spanish text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 51
$numRows  = $lastRow - $firstRow + 1
$numCols  = 16   # columns A..P

# Map: new row offset (0-based, 0 => row 2) -> old row number it must contain.
$srcRow = @(7,8,9,10,11,12,13,14,2,15,16,17,18,19,20,21,22,23,24,3,25,26,27,28,29,30,31,4,32,5,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,6)

# Snapshot the whole block before writing anything back (so we never read
# a cell we've already overwritten).
$srcRange = $ws.Range("A$firstRow`:P$lastRow")
$orig = $srcRange.Value2

# Build the permuted block in memory.
$new = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    $oldRowNum = $srcRow[$i]
    $oldOffset = $oldRowNum - $firstRow + 1   # 1-based row index into $orig
    for ($c = 1; $c -le $numCols; $c++) {
        $new[$i, $c - 1] = $orig[$oldOffset, $c]
    }
}

# Columns A, B, D, E hold id/date-like text (e.g. "1497", "4/4/2024", "2",
# "784804268") that Excel would otherwise auto-coerce into numbers/dates.
# Force them to remain plain text, matching the source file's inlineStr
# cells, before writing the values in.
$ws.Range("A$firstRow`:B$lastRow").NumberFormat = "@"
$ws.Range("D$firstRow`:E$lastRow").NumberFormat = "@"

$dstRange = $ws.Range("A$firstRow`:P$lastRow")
$dstRange.Value2 = $new
